# Apply updated TPM-based NATMI ligand-receptor statistics to Gnai2-C5ar1 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
    ,@("ECs","Gnai2","C5ar1","ECs",3,1,169.915657,509.746971,0.4441184931734509,0.4441184931734509,3,1,0.3893403333333333,1.168021,0.01344665141573655,0.01344665141573655,66.15501853493234,595.395166814391,0.005971906564985569,0.00597190656498557)
    ,@("ECs","Gnai2","C5ar1","FAPs",3,1,169.915657,509.746971,0.4441184931734509,0.4441184931734509,3,1,1.193104333333333,3.579313,0.04120625760907917,0.04120625760907917,202.7271066678803,1824.543960010923,0.01830046103866128,0.01830046103866129)
    ,@("ECs","Gnai2","C5ar1","MuSCs",3,1,169.915657,509.746971,0.4441184931734509,0.4441184931734509,1,0.3333333333333333,0.087271,0.261813,0.003014079496094877,0.003014079496094878,14.828709302047,133.458383718423,0.001338608444110651,0.001338608444110651)
    ,@("ECs","Gnai2","C5ar1","Resolving-Mac",3,1,169.915657,509.746971,0.4441184931734509,0.4441184931734509,3,1,27.28472966666666,81.85418899999999,0.9423330114790893,0.9423330114790894,4636.102767379058,41724.92490641151,0.4185075171256933,0.4185075171256935)
    ,@("FAPs","Gnai2","C5ar1","ECs",3,1,68.382243,205.146729,0.1787346690539575,0.1787346690539575,3,1,0.3893403333333333,1.168021,0.01344665141573655,0.01344665141573655,26.623965283701,239.615687553309,0.002403382790675603,0.002403382790675603)
    ,@("FAPs","Gnai2","C5ar1","FAPs",3,1,68.382243,205.146729,0.1787346690539575,0.1787346690539575,3,1,1.193104333333333,3.579313,0.04120625760907917,0.04120625760907917,81.58715044635301,734.284354017177,0.007364986816710883,0.007364986816710883)
    ,@("FAPs","Gnai2","C5ar1","MuSCs",3,1,68.382243,205.146729,0.1787346690539575,0.1787346690539575,1,0.3333333333333333,0.087271,0.261813,0.003014079496094877,0.003014079496094878,5.967786728853,53.710080559677,0.0005387205012368369,0.000538720501236837)
    ,@("FAPs","Gnai2","C5ar1","Resolving-Mac",3,1,68.382243,205.146729,0.1787346690539575,0.1787346690539575,3,1,27.28472966666666,81.85418899999999,0.9423330114790893,0.9423330114790894,1865.791014255309,16792.11912829778,0.1684275789453342,0.1684275789453342)
    ,@("MuSCs","Gnai2","C5ar1","ECs",3,1,53.27463399999999,159.823902,0.1392470275793777,0.1392470275793778,3,1,0.3893403333333333,1.168021,0.01344665141573655,0.01344665141573655,20.74196375977133,186.677673837942,0.001872406240537347,0.001872406240537347)
    ,@("MuSCs","Gnai2","C5ar1","FAPs",3,1,53.27463399999999,159.823902,0.1392470275793777,0.1392470275793778,3,1,1.193104333333333,3.579313,0.04120625760907917,0.04120625760907917,63.56219668214732,572.0597701393259,0.005737848889734391,0.005737848889734391)
    ,@("MuSCs","Gnai2","C5ar1","MuSCs",3,1,53.27463399999999,159.823902,0.1392470275793777,0.1392470275793778,1,0.3333333333333333,0.087271,0.261813,0.003014079496094877,0.003014079496094878,4.649330583814,41.843975254326,0.0004197016107191603,0.0004197016107191604)
    ,@("MuSCs","Gnai2","C5ar1","Resolving-Mac",3,1,53.27463399999999,159.823902,0.1392470275793777,0.1392470275793778,3,1,27.28472966666666,81.85418899999999,0.9423330114790893,0.9423330114790894,1453.583986780608,13082.25588102547,0.1312170708383868,0.1312170708383869)
    ,@("Resolving-Mac","Gnai2","C5ar1","ECs",3,1,91.01828266666666,273.054848,0.2378998101932138,0.2378998101932138,3,1,0.3893403333333333,1.168021,0.01344665141573655,0.01344665141573655,35.43708851286755,318.933796615808,0.003198955819538036,0.003198955819538036)
    ,@("Resolving-Mac","Gnai2","C5ar1","FAPs",3,1,91.01828266666666,273.054848,0.2378998101932138,0.2378998101932138,3,1,1.193104333333333,3.579313,0.04120625760907917,0.04120625760907917,108.5943074621582,977.3487671594239,0.009802960863972605,0.009802960863972606)
    ,@("Resolving-Mac","Gnai2","C5ar1","MuSCs",3,1,91.01828266666666,273.054848,0.2378998101932138,0.2378998101932138,1,0.3333333333333333,0.087271,0.261813,0.003014079496094877,0.003014079496094878,7.943256546602667,71.489308919424,0.0007170489400282287,0.0007170489400282289)
    ,@("Resolving-Mac","Gnai2","C5ar1","Resolving-Mac",3,1,91.01828266666666,273.054848,0.2378998101932138,0.2378998101932138,3,1,27.28472966666666,81.85418899999999,0.9423330114790893,0.9423330114790894,2483.409237284252,22350.68313555827,0.2241808445696749,0.224180844569675)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $addr = $cols[$c] + $rowNum
        $ws.Range($addr).Value = $rowData[$c]
    }
}
